$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.494.47'
$ws.Range('E2').Value = '  -0.21%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.733.48'
$ws.Range('E3').Value = '  -0.33%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.09%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.84'
$ws.Range('E5').Value = '  +0.41%  '

# Row 6
$ws.Range('E6').Value = '  +0.05%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4874'
$ws.Range('E7').Value = '  +1.67%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2662'
$ws.Range('E8').Value = '  -0.79%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06211'

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.735.80'
$ws.Range('E10').Value = '  -0.18%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07042'
$ws.Range('E11').Value = '  -0.98%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.70'
$ws.Range('E12').Value = '  -0.89%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.602'
$ws.Range('E13').Value = '  +1.12%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6085'
$ws.Range('E14').Value = '  -1.90%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.27'

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('E16').Value = '  +0.00%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.488.75'
$ws.Range('E17').Value = '  -0.25%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.001'
$ws.Range('E18').Value = '  +0.10%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007235'
$ws.Range('E19').Value = '  +4.93%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.50'
$ws.Range('E20').Value = '  -2.37%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.959.12'
$ws.Range('E21').Value = '  -0.17%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.532'
$ws.Range('E22').Value = '  -1.58%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.755'
$ws.Range('E23').Value = '  -1.57%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.236'
$ws.Range('E24').Value = '  -2.25%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '139.44'
$ws.Range('E25').Value = '  +2.77%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.42'
$ws.Range('E26').Value = '  +0.20%  '

# Row 27
$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.414'
$ws.Range('E27').Value = '  -0.80%  '

# Row 28
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.778'
$ws.Range('E28').Value = '  -1.78%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '107.83'
$ws.Range('E29').Value = '  +0.70%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.977'
$ws.Range('E30').Value = '  -0.73%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08062'
$ws.Range('E31').Value = '  +2.54%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.681'
$ws.Range('E32').Value = '  -1.52%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04565'
$ws.Range('E33').Value = '  -0.30%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.000'
$ws.Range('E34').Value = '  +0.04%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.617'
$ws.Range('E35').Value = '  +0.11%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.008'
$ws.Range('E36').Value = '  +0.74%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6379'
$ws.Range('E37').Value = '  -0.10%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.9039'
$ws.Range('E38').Value = '  -3.25%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.034'
$ws.Range('E39').Value = '  +2.48%  '

# Row 40
$ws.Range('E40').Value = '  -1.20%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.004'
$ws.Range('E41').Value = '  +0.01%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01505'
$ws.Range('E42').Value = '  -0.37%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.47'
$ws.Range('E43').Value = '  -10.36%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.422'
$ws.Range('E44').Value = '  -5.77%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.3894'
$ws.Range('E45').Value = '  -0.48%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.937'
$ws.Range('E46').Value = '  +0.66%  '

# Row 47
$ws.Range('E47').Value = '  -1.81%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05391'
$ws.Range('E48').Value = '  +1.14%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.55'
$ws.Range('E49').Value = '  -0.61%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.804'
$ws.Range('E50').Value = '  -0.99%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.250'
$ws.Range('E51').Value = '  -0.89%  '
